$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the arithmetic instructions (rows 17-20) as Implemented ("Y") and
# record whether they've been Tested ("Y" = fully, "?" = unclear/partial).
$ws.Range("E17").Value = "Y"
$ws.Range("F17").Value = "?"

$ws.Range("E18").Value = "Y"
$ws.Range("F18").Value = "Y"

$ws.Range("E19").Value = "Y"
$ws.Range("F19").Value = "?"

$ws.Range("E20").Value = "Y"
$ws.Range("F20").Value = "Y"

# Reflect the author's last active selection when they finished editing.
$ws.Range("F20").Select()
